$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 192 (pushes the existing row 192.."last" down by one,
# e.g. old row 192 becomes row 193, ..., old row 285 becomes row 286).
$ws.Rows.Item(192).EntireRow.Insert()

# Populate the newly inserted row 192 with the new weekly price observation.
# Columns A-L, Q, R, T repeat the same catalog/metadata values as the row that
# used to sit at 192 (now shifted to 193); D, M, N, O, P, S carry the new figures.
$ws.Cells.Item(192, 1).Value = 5
$ws.Cells.Item(192, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(192, 3).Value = "Maule"
$ws.Cells.Item(192, 4).Value = 44813
$ws.Cells.Item(192, 5).Value = 7
$ws.Cells.Item(192, 6).Value = "Fruta"
$ws.Cells.Item(192, 7).Value = 100108
$ws.Cells.Item(192, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(192, 9).Value = 100108005
$ws.Cells.Item(192, 10).Value = "Piña"
$ws.Cells.Item(192, 11).Value = "Caramelo"
$ws.Cells.Item(192, 12).Value = "Segunda"
$ws.Cells.Item(192, 13).Value = 410
$ws.Cells.Item(192, 14).Value = 18000
$ws.Cells.Item(192, 15).Value = 19000
$ws.Cells.Item(192, 16).Value = 18634
$ws.Cells.Item(192, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(192, 18).Value = "Ecuador"
$ws.Cells.Item(192, 19).Value = 1331
$ws.Cells.Item(192, 20).Value = 14
